$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I5").Value = 1.289671547221794
$ws.Range("J5").Value = 0.5206953613327997
$ws.Range("K5").Value = 0.07372489015042043
$ws.Range("L5").Value = 2.19796163554654
